# Feria Lagunitas de Puerto Montt - Ciboulette
# A new weekly price observation was recorded, inserted as a new row 112
# (shifting all subsequent rows down by one, rows 112-134 -> 113-135).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 112; everything below (old rows 112-134)
# shifts down to 113-135, carrying its formatting (e.g. the date style on D).
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new observation.
$ws.Range("A112").Value = 4
$ws.Range("B112").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C112").Value = "Los Lagos"
$ws.Range("D112").Value = 44511
$ws.Range("E112").Value = 10
$ws.Range("F112").Value = 100112039
$ws.Range("G112").Value = "Ciboulette"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 80
$ws.Range("K112").Value = 2500
$ws.Range("L112").Value = 2500
$ws.Range("M112").Value = 2500
$ws.Range("N112").Value = "`$/docena de atados"
$ws.Range("O112").Value = "Región Metropolitana"
$ws.Range("P112").Value = 833
$ws.Range("Q112").Value = 3
$ws.Range("R112").Value = "Hortaliza"
